$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D15 (shared string 15) text
$ws.Range("D15").Value = "Created a spreadsheet to analyse baseline data, found a bug in simulator from the statistical data. Fixed simulator and generated another data set. Logging and reporting is finished for this iteration."

# Update C15 value from 1 to 2
$ws.Range("C15").Value = 2

# Update selection to D15
$ws.Range("D15").Select()

# Autofit row 15 height to reflect the now-wrapped longer text
$ws.Rows(15).RowHeight = 42.75
